# Apply Iteration 3 test-case result updates to "Login & Account Management" sheet
# (matches commit: "Update test result for iter 3")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2 G: 'Fail' -> 'Pass'
$ws.Range("G2").Value = 'Pass'

# Row 2 H: 'unknown' -> None
$ws.Range("H2").ClearContents() | Out-Null

# Row 6 G: 'Fail' -> 'Pass'
$ws.Range("G6").Value = 'Pass'

# Row 6 H: 'unknown' -> None
$ws.Range("H6").ClearContents() | Out-Null

# Row 7 F: 'Login Page Redirect To http://ec2-35-161-162-8.us-west-2.compute.amazonaws.com/U-Link/index.html' -> 'Login Page Redirect To http://ec2-35-164-141-199.us-west-2.compute.amazonaws.com/U-Link/index.html'
$ws.Range("F7").Value = 'Login Page Redirect To http://ec2-35-164-141-199.us-west-2.compute.amazonaws.com/U-Link/index.html'

# Row 8 F: '"username is taken"\n\nNote: This account appears in account mangement tab' -> '"Account successfully created" shown , redirect to accountmangement.html'
$ws.Range("F8").Value = '"Account successfully created" shown , redirect to accountmangement.html'

# Row 8 G: 'Fail' -> 'Pass'
$ws.Range("G8").Value = 'Pass'

# Row 8 H: 'UI' -> None
$ws.Range("H8").ClearContents() | Out-Null

# Row 9 F: 'After clicking login button, it stays at the same page' -> 'Incorrect username/password'
$ws.Range("F9").Value = 'Incorrect username/password'

# Row 9 G: 'Fail' -> 'Pass'
$ws.Range("G9").Value = 'Pass'

# Row 9 H: 'unknown' -> None
$ws.Range("H9").ClearContents() | Out-Null

# Row 10 F: 'Login Page Redirect To http://ec2-35-161-162-8.us-west-2.compute.amazonaws.com/U-Link/index.html' -> 'Login Page Redirect To http://ec2-35-164-141-199.us-west-2.compute.amazonaws.com/U-Link/index.html'
$ws.Range("F10").Value = 'Login Page Redirect To http://ec2-35-164-141-199.us-west-2.compute.amazonaws.com/U-Link/index.html'

# Row 12 G: 'Fail' -> 'Pass'
$ws.Range("G12").Value = 'Pass'

# Row 12 H: 'unknown' -> None
$ws.Range("H12").ClearContents() | Out-Null

# Row 13 G: 'Fail' -> 'pass'
$ws.Range("G13").Value = 'pass'

# Row 13 H: 'ecf' -> None
$ws.Range("H13").ClearContents() | Out-Null

# Row 14 G: 'Fail' -> 'pass'
$ws.Range("G14").Value = 'pass'

# Row 14 H: 'ecf' -> None
$ws.Range("H14").ClearContents() | Out-Null

# Row 15 E: '"Account successfully deleted." shown' -> 'Account disappear from the table '
$ws.Range("E15").Value = 'Account disappear from the table '

# Row 15 F: 'the account was removed without any success message' -> 'the account was removed '
$ws.Range("F15").Value = 'the account was removed '

# Row 15 G: 'Fail' -> 'Pass'
$ws.Range("G15").Value = 'Pass'

# Row 15 H: 'UI' -> None
$ws.Range("H15").ClearContents() | Out-Null

# Row 16 F: 'A blue color text was shown for 1 second under  "Create New Account" but was missing after that. Page did not redirect to any page.\n\nNote: This account appears in account mangement tab' -> '"Account successfully created" shown '
$ws.Range("F16").Value = '"Account successfully created" shown '

# Row 16 G: 'Fail' -> 'Pass'
$ws.Range("G16").Value = 'Pass'

# Row 16 H: 'UI' -> None
$ws.Range("H16").ClearContents() | Out-Null

# Row 17 F: 'A blue color text was shown for 1 second under  "Create New Account" but was missing after that. Page did not redirect to any page.\n\nNote: This account appears in account mangement tab' -> '"Account successfully created" shown '
$ws.Range("F17").Value = '"Account successfully created" shown '

# Row 17 G: 'Fail' -> 'Pass'
$ws.Range("G17").Value = 'Pass'

# Row 17 H: 'UI' -> None
$ws.Range("H17").ClearContents() | Out-Null

# Row 18 D: 'Login with Test Case 16\nEmail: anotheruser@hotmail.com\nPassword: password \nClick Login button' -> 'Login with Test Case 16\nEmail: anotheruser\nPassword: password \nClick Login button'
$ws.Range("D18").Value = 'Login with Test Case 16' + [char]10 + 'Email: anotheruser' + [char]10 + 'Password: password ' + [char]10 + 'Click Login button'

# Row 21 F: '"Username is taken" was shown for 1 second under  "Create New Account" but was missing after that. Page did not redirect to any page.\n ' -> 'Username is taken '
$ws.Range("F21").Value = 'Username is taken '

# Row 21 G: 'Fail' -> 'Pass'
$ws.Range("G21").Value = 'Pass'

# Row 21 H: 'UI' -> None
$ws.Range("H21").ClearContents() | Out-Null

# Row 22 F: '"Username is taken" was shown for 1 second under  "Create New Account" but was missing after that. Page did not redirect to any page.\n ' -> 'Username is taken '
$ws.Range("F22").Value = 'Username is taken '

# Row 22 G: 'Fail' -> 'Pass'
$ws.Range("G22").Value = 'Pass'

# Row 22 H: 'UI' -> None
$ws.Range("H22").ClearContents() | Out-Null

# Row 23 G: 'Fail' -> 'Pass'
$ws.Range("G23").Value = 'Pass'

# Row 23 H: 'unknown' -> None
$ws.Range("H23").ClearContents() | Out-Null

# Row 29 F: 'Nothing happens ' -> 'Password cannot be empty'
$ws.Range("F29").Value = 'Password cannot be empty'

# Row 29 G: 'Fail' -> 'Pass'
$ws.Range("G29").Value = 'Pass'

# Update the row heights for rows whose "Actual Output" text got shorter so the
# sheet re-wraps/auto-fits like Excel would after a content edit.
$ws.Rows("16:17").AutoFit() | Out-Null
$ws.Rows("21:22").AutoFit() | Out-Null

# Restore the view: scroll so row 15 / column A is the top-left visible cell,
# and select E17 (matches the saved selection in the edited workbook).
$ws.Activate()
$ws.Range("E17").Select()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
